# Add a new column J to the "Suivi" sheet:
#  - J1 header gets a new timestamp, with the same header style as the
#    other header cells (column I).
#  - J2:J100 get a price snapshot that duplicates the current column I
#    (latest known price) values.
#  - J101:J204 stay blank, matching the existing blank cells in column I
#    for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 100
$lastRow = 204

# Header cell: copy formatting from I1 (so it keeps style index 1), then
# overwrite with the new timestamp text.
$ws.Range("I1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "2026-01-27 23:11:44"

# Price rows: duplicate column I's numeric values into column J.
for ($r = 2; $r -le $lastDataRow; $r++) {
    $price = $ws.Cells.Item($r, 9).Value2
    $ws.Cells.Item($r, 10).Value2 = $price
}

# Trailing rows: column I is blank for these rows, so mirror that blank
# cell into column J (keeps the sheet's used range consistent).
for ($r = ($lastDataRow + 1); $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Copy($ws.Cells.Item($r, 10))
}
